$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Cells.Item(3, 1).Value = 111782750
$ws.Cells.Item(3, 2).Value = 82949
$ws.Cells.Item(3, 4).Value = 'NT'
$ws.Cells.Item(3, 5).Value = 5589
$ws.Cells.Item(3, 6).Value = 'Rödbrun klubbdyna'
$ws.Cells.Item(3, 7).Value = 'Trichoderma nybergianum'
$ws.Cells.Item(3, 8).Value = '(T.Ulvinen & H.L.Chamb.) Jaklitsch & Voglmayr'
$ws.Cells.Item(3, 10).Value = ""
$ws.Cells.Item(3, 14).Value = ""
$ws.Cells.Item(3, 17).Value = 505007.618534557
$ws.Cells.Item(3, 18).Value = 7018756.52538473
$ws.Cells.Item(3, 32).Value = ""

# Row 4
$ws.Cells.Item(4, 1).Value = 111783071
$ws.Cells.Item(4, 2).Value = 90300
$ws.Cells.Item(4, 5).Value = 4745
$ws.Cells.Item(4, 6).Value = 'Tallriska'
$ws.Cells.Item(4, 7).Value = 'Lactarius musteus'
$ws.Cells.Item(4, 8).Value = 'Fr.'
$ws.Cells.Item(4, 10).Value = ""
$ws.Cells.Item(4, 14).Value = ""
$ws.Cells.Item(4, 17).Value = 505060.2648977584
$ws.Cells.Item(4, 18).Value = 7018787.191973396
$ws.Cells.Item(4, 32).Value = ""

# Row 5
$ws.Cells.Item(5, 1).Value = 111783769
$ws.Cells.Item(5, 2).Value = 90665
$ws.Cells.Item(5, 4).Value = 'VU'
$ws.Cells.Item(5, 5).Value = 1435
$ws.Cells.Item(5, 6).Value = 'Bitter taggsvamp'
$ws.Cells.Item(5, 7).Value = 'Hydnellum fennicum'
$ws.Cells.Item(5, 8).Value = '(P.Karst.) E.Larss., K.H.Larss. & Kõljalg'
$ws.Cells.Item(5, 17).Value = 505036.7939151306
$ws.Cells.Item(5, 18).Value = 7018819.987804689

# Row 6
$ws.Cells.Item(6, 1).Value = 111782876
$ws.Cells.Item(6, 2).Value = 90651
$ws.Cells.Item(6, 5).Value = 1968
$ws.Cells.Item(6, 6).Value = 'Grantaggsvamp'
$ws.Cells.Item(6, 7).Value = 'Bankera violascens'
$ws.Cells.Item(6, 8).Value = '(Alb. & Schwein. : Fr.) Pouzar'
$ws.Cells.Item(6, 17).Value = 505022.9813479512
$ws.Cells.Item(6, 18).Value = 7018724.615566149

# Row 7
$ws.Cells.Item(7, 1).Value = 111814356
$ws.Cells.Item(7, 2).Value = 90678
$ws.Cells.Item(7, 4).Value = 'LC'
$ws.Cells.Item(7, 5).Value = 4366
$ws.Cells.Item(7, 6).Value = 'Skarp dropptaggsvamp'
$ws.Cells.Item(7, 7).Value = 'Hydnellum peckii'
$ws.Cells.Item(7, 8).Value = 'Banker'
$ws.Cells.Item(7, 17).Value = 505204.4099656619
$ws.Cells.Item(7, 18).Value = 7018772.129998797

# Row 8
$ws.Cells.Item(8, 1).Value = 111814350
$ws.Cells.Item(8, 2).Value = 89401
$ws.Cells.Item(8, 4).Value = 'NT'
$ws.Cells.Item(8, 5).Value = 1108
$ws.Cells.Item(8, 6).Value = 'Harticka'
$ws.Cells.Item(8, 7).Value = 'Pelloporus leporinus'
$ws.Cells.Item(8, 8).Value = '(Fr.) Krieglst.'
$ws.Cells.Item(8, 17).Value = 504984.8875472886
$ws.Cells.Item(8, 18).Value = 7018893.217038274

# Row 9
$ws.Cells.Item(9, 1).Value = 111814417
$ws.Cells.Item(9, 2).Value = 103288
$ws.Cells.Item(9, 5).Value = 221144
$ws.Cells.Item(9, 6).Value = 'Grönpyrola'
$ws.Cells.Item(9, 7).Value = 'Pyrola chlorantha'
$ws.Cells.Item(9, 8).Value = 'Sw.'
$ws.Cells.Item(9, 17).Value = 505110.448201828
$ws.Cells.Item(9, 18).Value = 7018666.73204405

# Row 10
$ws.Cells.Item(10, 1).Value = 111814375
$ws.Cells.Item(10, 2).Value = 96370
$ws.Cells.Item(10, 5).Value = 219847
$ws.Cells.Item(10, 6).Value = 'Tvåblad'
$ws.Cells.Item(10, 7).Value = 'Neottia ovata'
$ws.Cells.Item(10, 8).Value = '(L.) Buff. & Fingerh.'
$ws.Cells.Item(10, 17).Value = 504999.9977373667
$ws.Cells.Item(10, 18).Value = 7018733.575208749

# Row 12
$ws.Cells.Item(12, 1).Value = 111814432
$ws.Cells.Item(12, 2).Value = 95674
$ws.Cells.Item(12, 4).Value = 'LC'
$ws.Cells.Item(12, 5).Value = 222741
$ws.Cells.Item(12, 6).Value = 'Finbräken'
$ws.Cells.Item(12, 7).Value = 'Cystopteris montana'
$ws.Cells.Item(12, 8).Value = '(Lam.) Desv.'
$ws.Cells.Item(12, 17).Value = 505015.75484597
$ws.Cells.Item(12, 18).Value = 7018737.647699019

# Row 13
$ws.Cells.Item(13, 1).Value = 111814395
$ws.Cells.Item(13, 2).Value = 103288
$ws.Cells.Item(13, 5).Value = 221144
$ws.Cells.Item(13, 6).Value = 'Grönpyrola'
$ws.Cells.Item(13, 7).Value = 'Pyrola chlorantha'
$ws.Cells.Item(13, 8).Value = 'Sw.'

# Row 14
$ws.Cells.Item(14, 1).Value = 111814388
$ws.Cells.Item(14, 2).Value = 78578
$ws.Cells.Item(14, 4).Value = 'NT'
$ws.Cells.Item(14, 5).Value = 6458
$ws.Cells.Item(14, 6).Value = 'Lunglav'
$ws.Cells.Item(14, 7).Value = 'Lobaria pulmonaria'
$ws.Cells.Item(14, 8).Value = '(L.) Hoffm.'
$ws.Cells.Item(14, 17).Value = 505182.7410700406
$ws.Cells.Item(14, 18).Value = 7018803.578552675
$ws.Cells.Item(14, 36).Value = 'sälg'
$ws.Cells.Item(14, 37).Value = 'Salix caprea'
$ws.Cells.Item(14, 41).Value = 'Salix caprea'

# Row 16
$ws.Cells.Item(16, 1).Value = 111814415
$ws.Cells.Item(16, 2).Value = 56543
$ws.Cells.Item(16, 4).Value = 'NT'
$ws.Cells.Item(16, 5).Value = 103021
$ws.Cells.Item(16, 6).Value = 'Talltita'
$ws.Cells.Item(16, 7).Value = 'Poecile montanus'
$ws.Cells.Item(16, 8).Value = '(Conrad von Baldenstein, 1827)'
$ws.Cells.Item(16, 17).Value = 505100.2035408606
$ws.Cells.Item(16, 18).Value = 7018878.55609256

# Row 17
$ws.Cells.Item(17, 1).Value = 111814359
$ws.Cells.Item(17, 2).Value = 90300
$ws.Cells.Item(17, 5).Value = 4745
$ws.Cells.Item(17, 6).Value = 'Tallriska'
$ws.Cells.Item(17, 7).Value = 'Lactarius musteus'
$ws.Cells.Item(17, 8).Value = 'Fr.'
$ws.Cells.Item(17, 17).Value = 505073.4975346876
$ws.Cells.Item(17, 18).Value = 7018678.36784017
$ws.Cells.Item(17, 36).Value = ""
$ws.Cells.Item(17, 37).Value = ""
$ws.Cells.Item(17, 41).Value = ""

# Row 18
$ws.Cells.Item(18, 1).Value = 111814428
$ws.Cells.Item(18, 2).Value = 90666
$ws.Cells.Item(18, 5).Value = 4364
$ws.Cells.Item(18, 6).Value = 'Dropptaggsvamp'
$ws.Cells.Item(18, 7).Value = 'Hydnellum ferrugineum'
$ws.Cells.Item(18, 8).Value = '(Fr.:Fr.) P. Karst.'
$ws.Cells.Item(18, 17).Value = 504987.8220338543
$ws.Cells.Item(18, 18).Value = 7018743.451279385

# Row 19
$ws.Cells.Item(19, 1).Value = 111814362
$ws.Cells.Item(19, 2).Value = 103288
$ws.Cells.Item(19, 5).Value = 221144
$ws.Cells.Item(19, 6).Value = 'Grönpyrola'
$ws.Cells.Item(19, 7).Value = 'Pyrola chlorantha'
$ws.Cells.Item(19, 8).Value = 'Sw.'
$ws.Cells.Item(19, 17).Value = 504958.3523041067
$ws.Cells.Item(19, 18).Value = 7018869.788911887

# Row 20
$ws.Cells.Item(20, 1).Value = 111814434
$ws.Cells.Item(20, 2).Value = 96253
$ws.Cells.Item(20, 5).Value = 504
$ws.Cells.Item(20, 6).Value = 'Guckusko'
$ws.Cells.Item(20, 7).Value = 'Cypripedium calceolus'
$ws.Cells.Item(20, 8).Value = 'L.'
$ws.Cells.Item(20, 17).Value = 505012.6056710624
$ws.Cells.Item(20, 18).Value = 7018735.393927739

# Row 21
$ws.Cells.Item(21, 1).Value = 111814351
$ws.Cells.Item(21, 2).Value = 89845
$ws.Cells.Item(21, 4).Value = 'VU'
$ws.Cells.Item(21, 5).Value = 1209
$ws.Cells.Item(21, 6).Value = 'Rynkskinn'
$ws.Cells.Item(21, 7).Value = 'Phlebia centrifuga'
$ws.Cells.Item(21, 8).Value = 'P.Karst.'
$ws.Cells.Item(21, 17).Value = 504984.8875472886
$ws.Cells.Item(21, 18).Value = 7018893.217038274

# Row 22
$ws.Cells.Item(22, 1).Value = 111814402
$ws.Cells.Item(22, 2).Value = 90651
$ws.Cells.Item(22, 5).Value = 1968
$ws.Cells.Item(22, 6).Value = 'Grantaggsvamp'
$ws.Cells.Item(22, 7).Value = 'Bankera violascens'
$ws.Cells.Item(22, 8).Value = '(Alb. & Schwein. : Fr.) Pouzar'
$ws.Cells.Item(22, 17).Value = 505200.3682009591
$ws.Cells.Item(22, 18).Value = 7018764.927175661

# Row 23
$ws.Cells.Item(23, 1).Value = 111814411
$ws.Cells.Item(23, 2).Value = 96253
$ws.Cells.Item(23, 4).Value = 'LC'
$ws.Cells.Item(23, 5).Value = 504
$ws.Cells.Item(23, 6).Value = 'Guckusko'
$ws.Cells.Item(23, 7).Value = 'Cypripedium calceolus'
$ws.Cells.Item(23, 8).Value = 'L.'
$ws.Cells.Item(23, 17).Value = 505014.8575873387
$ws.Cells.Item(23, 18).Value = 7018735.397438973

# Row 24
$ws.Cells.Item(24, 1).Value = 111814348
$ws.Cells.Item(24, 2).Value = 96370
$ws.Cells.Item(24, 5).Value = 219847
$ws.Cells.Item(24, 6).Value = 'Tvåblad'
$ws.Cells.Item(24, 7).Value = 'Neottia ovata'
$ws.Cells.Item(24, 8).Value = '(L.) Buff. & Fingerh.'
$ws.Cells.Item(24, 17).Value = 504944.9568800884
$ws.Cells.Item(24, 18).Value = 7018794.658574538
